$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.107.23"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3
$ws.Range("D3").Value = "1.789.12"

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'226.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "

# Row 6
$ws.Range("E6").Value = "  -0.55%  "

# Row 7
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").Value = "'32.27"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.04%  "

# Row 9
$ws.Range("E9").Value = "  +3.32%  "

# Row 10
$ws.Range("E10").Value = "  -2.83%  "

# Row 11
$ws.Range("E11").Value = "  +0.76%  "

# Row 12
$ws.Range("D12").Value = "2.047.73"
$ws.Range("E12").Value = "  -0.14%  "

# Row 13
$ws.Range("D13").Value = "'11.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.21%  "

# Row 14
$ws.Range("D14").Value = "1.785.54"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15
$ws.Range("D15").Value = "34.089.29"

# Row 16
$ws.Range("D16").Value = "'0.621"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "

# Row 17
$ws.Range("E17").Value = "  +0.08%  "

# Row 18
$ws.Range("D18").Value = "'67.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "

# Row 19
$ws.Range("D19").Value = "'243.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "

# Row 20
$ws.Range("E20").Value = "  -0.91%  "

# Row 21
$ws.Range("D21").Value = "'10.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.79%  "

# Row 22
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("E24").Value = "  -3.11%  "

# Row 25
$ws.Range("D25").Value = "'160.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.58%  "

# Row 26
$ws.Range("D26").Value = "'7.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.56%  "

# Row 27
$ws.Range("D27").Value = "'16.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "

# Row 28
$ws.Range("E28").Value = "  +0.85%  "

# Row 29
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("D30").Value = "'1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.05%  "

# Row 31
$ws.Range("E31").Value = "  -0.35%  "

# Row 32
$ws.Range("E32").Value = "  -0.35%  "

# Row 33
$ws.Range("D33").Value = "'3.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.17%  "

# Row 34
$ws.Range("E34").Value = "  +0.51%  "

# Row 35
$ws.Range("D35").Value = "1.405.05"
$ws.Range("E35").Value = "  +0.75%  "

# Row 36
$ws.Range("E36").Value = "  +0.99%  "

# Row 37
$ws.Range("E37").Value = "  -0.55%  "

# Row 38
$ws.Range("E38").Value = "  +1.54%  "

# Row 39
$ws.Range("E39").Value = "  +5.69%  "

# Row 40
$ws.Range("E40").Value = "  +1.31%  "

# Row 41
$ws.Range("D41").Value = "'79.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "

# Row 42
$ws.Range("D42").Value = "'0.921"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.69%  "

# Row 43
$ws.Range("E43").Value = "  +0.14%  "

# Row 44
$ws.Range("D44").Value = "'13.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.54%  "

# Row 45
$ws.Range("D45").Value = "0.0₆0139"
$ws.Range("E45").Value = "  -4.42%  "

# Row 46
$ws.Range("E46").Value = "  +2.68%  "

# Row 47
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "'0.0507"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.51%  "

# Row 48
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.12%  "

# Row 49
$ws.Range("D49").Value = "'106.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.20%  "

# Row 50
$ws.Range("D50").Value = "1.949.13"
$ws.Range("E50").Value = "  -0.25%  "

# Row 51
$ws.Range("E51").Value = "  +0.16%  "
